$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above current row 4, shifting rows 4-15 down to 5-16
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the new data point
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44473
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 100112013
$ws.Cells.Item(4, 7).Value = "Alcachofa"
$ws.Cells.Item(4, 8).Value = "Madrigal"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 160
$ws.Cells.Item(4, 11).Value = 11000
$ws.Cells.Item(4, 12).Value = 12000
$ws.Cells.Item(4, 13).Value = 11500
$ws.Cells.Item(4, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(4, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(4, 16).Value = 288
$ws.Cells.Item(4, 17).Value = 40
$ws.Cells.Item(4, 18).Value = "Hortaliza"
